# Weekly refresh of the daily "Fruta / hortaliza" price table.
# The data rows (2-19) get reshuffled: each row's Fecha (D) plus the
# Volumen / Precio / Unidad / Origen / Precio-Kg / Kg-unidad block
# (M:T) is replaced by the block that used to live on a different row,
# per the mapping below (destination row -> source row, both referring
# to the ORIGINAL / pre-edit layout).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    2  = 18
    3  = 4
    4  = 7
    5  = 19
    6  = 10
    7  = 14
    8  = 2
    9  = 12
    10 = 8
    11 = 9
    12 = 3
    13 = 17
    14 = 11
    15 = 16
    16 = 6
    17 = 5
    18 = 13
    19 = 15
}

# Snapshot the original D and M:T values for every data row before any
# writes happen, since several rows read from each other.
# Value2 (not Value) is used so numbers come back as plain doubles and
# the D-column date serial isn't reinterpreted/reformatted.
$snapshot = @{}
for ($r = 2; $r -le 19; $r++) {
    $row = @{}
    $row['D'] = $ws.Cells.Item($r, 4).Value2
    for ($c = 13; $c -le 20; $c++) {
        $row[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $row
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $src = $snapshot[$srcRow]

    $ws.Cells.Item($destRow, 4).Value2 = $src['D']
    for ($c = 13; $c -le 20; $c++) {
        $ws.Cells.Item($destRow, $c).Value2 = $src[$c]
    }
}
